$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Slide 1 - title slide: bump the presenter/date line by a day
# -----------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$dateShape = $s1.Shapes.Item(5)
$dateTr = $dateShape.TextFrame.TextRange
$dateLen = $dateTr.Text.Length
$dateTr.Characters(1, $dateLen).Text = '[Presenter Name] | November 18, 2025'

# -----------------------------------------------------------------
# 2) Slide 9 - Investment Summary table: add cost-category rows and
#    fill in the TOTAL INVESTMENT figures.
# -----------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tblShape = $s9.Shapes.Item(3)
$tbl = $tblShape.Table

# Resize the columns (values are points = EMU / 12700)
$tbl.Columns.Item(1).Width = 1567967 / 12700
$tbl.Columns.Item(2).Width = 1045311 / 12700
$tbl.Columns.Item(3).Width = 2090623 / 12700
$tbl.Columns.Item(4).Width = 1045311 / 12700
$tbl.Columns.Item(5).Width = 871093 / 12700
$tbl.Columns.Item(6).Width = 871093 / 12700
$tbl.Columns.Item(7).Width = 1219530 / 12700

# New cost-category rows, inserted (in order) right before the
# existing "TOTAL INVESTMENT" row (currently row 2).
$newRows = @(
    @('Professional Services', '$82,250', '($10,000)', '$72,250', '$0', '$0', '$72,250'),
    @('Cloud Infrastructure', '$18,528', '($3,690)', '$14,838', '$18,528', '$18,528', '$51,894'),
    @('Software Licenses', '$2,904', '$0', '$2,904', '$2,904', '$2,904', '$8,712'),
    @('Support & Maintenance', '$2,676', '$0', '$2,676', '$2,676', '$2,676', '$8,028')
)

$insertAt = 2
foreach ($rowValues in $newRows) {
    $tbl.Rows.Add($insertAt) | Out-Null
    for ($c = 1; $c -le 7; $c++) {
        $cell = $tbl.Cell($insertAt, $c)
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Text = $rowValues[$c - 1]
        $tr.Font.Size = 11
        $tr.Font.Bold = 0
    }
    $insertAt = $insertAt + 1
}

# Fill in the (previously blank) TOTAL INVESTMENT figures - this row
# has now been pushed down to row 6. (First column, the "TOTAL
# INVESTMENT" label, is already correct and left untouched.)
$totalRow = 6
$totalValues = @('$106,358', '($13,690)', '$92,668', '$24,108', '$24,108', '$140,884')
for ($c = 2; $c -le 7; $c++) {
    $tbl.Cell($totalRow, $c).Shape.TextFrame.TextRange.Text = $totalValues[$c - 2]
}
